# Loan RBI, Variable Instalments
#
# The "Repayment Schedule" sheet gets a new (blank) column inserted
# immediately before what used to be column N ("Late"), pushing the
# "Late" column to O and the "Outstanding" column (old P) to Q. The new
# column inherits the column width/formatting of the column to its left
# (M). The "Repayment Schedule" sheet also becomes the active sheet
# (it was "Transactions" before), with cell R7 selected there while the
# "Transactions" sheet keeps its own prior selection (B3) but is no
# longer the active tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Remember the width of column M so the freshly inserted column can
# pick up the same formatting (mirrors Excel's default "insert column"
# behaviour of carrying over the format of the column to the left).
$mColumnWidth = $ws.Columns("M:M").ColumnWidth

# Insert a new blank column before N; this shifts the old N ("Late")
# to O and the old P ("Outstanding") to Q.
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $mColumnWidth

# Make "Repayment Schedule" the active sheet/tab and select R7 there.
$ws.Activate()
$ws.Range("R7").Select()
